$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B:E stay as text when assigned, even for numeric-looking
# strings (prices/volumes are stored as text in the source data), without
# leaving any lasting NumberFormat/style change on the cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "62.606.07"
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("D3").Value = "3.013.48"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "594.84"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "151.98"
$ws.Range("E6").Value = "  +5.93%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "3.003.25"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").Value = "6.83"
$ws.Range("E10").Value = "  +13.09%  "
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "35.40"
$ws.Range("E14").Value = "  +3.93%  "
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "3.518.06"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "62.718.91"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "7.04"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "3.017.20"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "446.78"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "14.16"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").Value = "0.693"
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").Value = "7.49"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("D24").Value = "82.62"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("D25").Value = "11.18"
$ws.Range("E25").Value = "  +4.32%  "
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").Value = "  +3.86%  "
$ws.Range("D27").Value = "12.26"
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "7.48"
$ws.Range("E29").Value = "  +4.54%  "
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +10.39%  "
$ws.Range("D31").Value = "2.69"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "27.54"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "0.0₃0862"
$ws.Range("E35").Value = "  +6.45%  "
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("D37").Value = "5.85"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").Value = "3.09"
$ws.Range("E38").Value = "  +9.28%  "
$ws.Range("D39").Value = "2.09"
$ws.Range("E39").Value = "  +4.26%  "
$ws.Range("D40").Value = "0.129"
$ws.Range("E40").Value = "  +4.54%  "
$ws.Range("D41").Value = "50.37"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").Value = "8.98"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").Value = "44.57"
$ws.Range("E43").Value = "  +13.10%  "
$ws.Range("D44").Value = "0.304"
$ws.Range("E44").Value = "  +12.43%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0359"
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "386.80"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "2.700.31"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").Value = "133.37"
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("D49").Value = "26.42"
$ws.Range("E49").Value = "  +13.31%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "2.26"
$ws.Range("E51").Value = "  +5.57%  "

$dataRange.ClearFormats()
